# Scheduled-runner update: refresh cached Universalis market-price columns
# (currentAveragePrice / NQ / HQ, LevePrice NQ/HQ, LeveProfit NQ/HQ) across
# several Leve worksheets. Values are static snapshots (no formulas in the
# source cells), so we just overwrite the affected cells per-sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(11, 8).Value = 132.2  # H11: 15 -> 132.2
$ws.Cells.Item(11, 9).Value = 132.2  # I11: 15 -> 132.2
$ws.Cells.Item(11, 11).Value = 132.2  # K11: 15 -> 132.2
$ws.Cells.Item(11, 13).Value = 7.800000000000011  # M11: 125 -> 7.800000000000011

$ws.Cells.Item(17, 8).Value = 3571.4285  # H17: 3409 -> 3571.4285
$ws.Cells.Item(17, 10).Value = 3571.4285  # J17: 3409 -> 3571.4285
$ws.Cells.Item(17, 12).Value = 10714.2855  # L17: 10227 -> 10714.2855
$ws.Cells.Item(17, 14).Value = -11050.2855  # N17: -10563 -> -11050.2855

$ws.Cells.Item(21, 8).Value = 29921  # H21: 29920.5 -> 29921
$ws.Cells.Item(21, 10).Value = 29921  # J21: 29920.5 -> 29921
$ws.Cells.Item(21, 12).Value = 29921  # L21: 29920.5 -> 29921
$ws.Cells.Item(21, 14).Value = -30857  # N21: -30856.5 -> -30857

$ws.Cells.Item(23, 8).Value = 29921  # H23: 29920.5 -> 29921
$ws.Cells.Item(23, 10).Value = 29921  # J23: 29920.5 -> 29921
$ws.Cells.Item(23, 12).Value = 29921  # L23: 29920.5 -> 29921
$ws.Cells.Item(23, 14).Value = -30389  # N23: -30388.5 -> -30389

$ws.Cells.Item(29, 8).Value = 2083.3333  # H29: 1812.5 -> 2083.3333
$ws.Cells.Item(29, 10).Value = 0  # J29: 1000 -> 0
$ws.Cells.Item(29, 12).Value = 0  # L29: 3000 -> 0
$ws.Cells.Item(29, 14).Value = ""  # N29: clear (was -3562)

$ws.Cells.Item(32, 8).Value = 1235.4  # H32: 925 -> 1235.4
$ws.Cells.Item(32, 9).Value = 1000  # I32: 650 -> 1000
$ws.Cells.Item(32, 10).Value = 1294.25  # J32: 1200 -> 1294.25
$ws.Cells.Item(32, 11).Value = 1000  # K32: 650 -> 1000
$ws.Cells.Item(32, 12).Value = 1294.25  # L32: 1200 -> 1294.25
$ws.Cells.Item(32, 13).Value = -674  # M32: -324 -> -674
$ws.Cells.Item(32, 14).Value = -1946.25  # N32: -1852 -> -1946.25

$ws.Cells.Item(38, 8).Value = 9480  # H38: 12225 -> 9480
$ws.Cells.Item(38, 10).Value = 12498.333  # J38: 18125 -> 12498.333
$ws.Cells.Item(38, 12).Value = 37494.999  # L38: 54375 -> 37494.999
$ws.Cells.Item(38, 14).Value = -38238.999  # N38: -55119 -> -38238.999

$ws.Cells.Item(40, 8).Value = 8500  # H40: 8400 -> 8500
$ws.Cells.Item(40, 10).Value = 8500  # J40: 8400 -> 8500
$ws.Cells.Item(40, 12).Value = 8500  # L40: 8400 -> 8500
$ws.Cells.Item(40, 14).Value = -8850  # N40: -8750 -> -8850

$ws.Cells.Item(64, 8).Value = 6737.8  # H64: 6598.778 -> 6737.8
$ws.Cells.Item(64, 9).Value = 5929.3335  # I64: 4899.5 -> 5929.3335
$ws.Cells.Item(64, 11).Value = 5929.3335  # K64: 4899.5 -> 5929.3335
$ws.Cells.Item(64, 13).Value = -5681.3335  # M64: -4651.5 -> -5681.3335

$ws.Cells.Item(67, 8).Value = 6737.8  # H67: 6598.778 -> 6737.8
$ws.Cells.Item(67, 9).Value = 5929.3335  # I67: 4899.5 -> 5929.3335
$ws.Cells.Item(67, 11).Value = 5929.3335  # K67: 4899.5 -> 5929.3335
$ws.Cells.Item(67, 13).Value = -5071.3335  # M67: -4041.5 -> -5071.3335

$ws.Cells.Item(92, 8).Value = 1210.4667  # H92: 1167.625 -> 1210.4667
$ws.Cells.Item(92, 9).Value = 1116  # I92: 1070.5385 -> 1116
$ws.Cells.Item(92, 11).Value = 1116  # K92: 1070.5385 -> 1116
$ws.Cells.Item(92, 13).Value = 132  # M92: 177.4614999999999 -> 132

$ws.Cells.Item(127, 8).Value = 2202  # H127: 1999 -> 2202
$ws.Cells.Item(127, 10).Value = 4106.5  # J127: 4996 -> 4106.5
$ws.Cells.Item(127, 12).Value = 12319.5  # L127: 14988 -> 12319.5
$ws.Cells.Item(127, 14).Value = -22239.5  # N127: -24908 -> -22239.5

$ws.Cells.Item(131, 8).Value = 6832.3335  # H131: 8931.333000000001 -> 6832.3335
$ws.Cells.Item(131, 9).Value = 4498.5  # I131: 7717.6 -> 4498.5
$ws.Cells.Item(131, 10).Value = 11500  # J131: 15000 -> 11500
$ws.Cells.Item(131, 11).Value = 13495.5  # K131: 23152.8 -> 13495.5
$ws.Cells.Item(131, 12).Value = 34500  # L131: 45000 -> 34500
$ws.Cells.Item(131, 13).Value = -8455.5  # M131: -18112.8 -> -8455.5
$ws.Cells.Item(131, 14).Value = -44580  # N131: -55080 -> -44580

$ws.Cells.Item(138, 8).Value = 3346.9644  # H138: 3361.4375 -> 3346.9644
$ws.Cells.Item(138, 9).Value = 2418.4  # I138: 2765.3333 -> 2418.4
$ws.Cells.Item(138, 10).Value = 3548.8262  # J138: 3499 -> 3548.8262
$ws.Cells.Item(138, 11).Value = 7255.200000000001  # K138: 8295.999899999999 -> 7255.200000000001
$ws.Cells.Item(138, 12).Value = 10646.4786  # L138: 10497 -> 10646.4786
$ws.Cells.Item(138, 13).Value = -2115.200000000001  # M138: -3155.999899999999 -> -2115.200000000001
$ws.Cells.Item(138, 14).Value = -20926.4786  # N138: -20777 -> -20926.4786

$ws.Cells.Item(141, 8).Value = 1998.5  # H141: 1873.5 -> 1998.5
$ws.Cells.Item(141, 9).Value = 1999  # I141: 1873.5 -> 1999
$ws.Cells.Item(141, 10).Value = 1998  # J141: 0 -> 1998
$ws.Cells.Item(141, 11).Value = 5997  # K141: 5620.5 -> 5997
$ws.Cells.Item(141, 12).Value = 5994  # L141: 0 -> 5994
$ws.Cells.Item(141, 13).Value = -817  # M141: -440.5 -> -817
$ws.Cells.Item(141, 14).Value = -16354  # N141: None -> -16354

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2915.946  # H32: 3001.6943 -> 2915.946
$ws.Cells.Item(32, 9).Value = 2915.946  # I32: 3001.6943 -> 2915.946
$ws.Cells.Item(32, 11).Value = 2915.946  # K32: 3001.6943 -> 2915.946
$ws.Cells.Item(32, 13).Value = -2628.946  # M32: -2714.6943 -> -2628.946

$ws.Cells.Item(74, 8).Value = 1464  # H74: 1522 -> 1464
$ws.Cells.Item(74, 9).Value = 1407.75  # I74: 1466 -> 1407.75
$ws.Cells.Item(74, 11).Value = 1407.75  # K74: 1466 -> 1407.75
$ws.Cells.Item(74, 13).Value = -533.75  # M74: -592 -> -533.75

$ws.Cells.Item(77, 8).Value = 1464  # H77: 1522 -> 1464
$ws.Cells.Item(77, 9).Value = 1407.75  # I77: 1466 -> 1407.75
$ws.Cells.Item(77, 11).Value = 7038.75  # K77: 7330 -> 7038.75
$ws.Cells.Item(77, 13).Value = -2670.75  # M77: -2962 -> -2670.75

$ws.Cells.Item(97, 8).Value = 4807.625  # H97: 5223 -> 4807.625
$ws.Cells.Item(97, 9).Value = 4494.4287  # I97: 4926.8335 -> 4494.4287
$ws.Cells.Item(97, 11).Value = 4494.4287  # K97: 4926.8335 -> 4494.4287
$ws.Cells.Item(97, 13).Value = -3998.4287  # M97: -4430.8335 -> -3998.4287

$ws.Cells.Item(110, 8).Value = 991.25  # H110: 648.4666999999999 -> 991.25
$ws.Cells.Item(110, 9).Value = 1199.5  # I110: 731.9167 -> 1199.5
$ws.Cells.Item(110, 10).Value = 366.5  # J110: 314.66666 -> 366.5
$ws.Cells.Item(110, 11).Value = 1199.5  # K110: 731.9167 -> 1199.5
$ws.Cells.Item(110, 12).Value = 366.5  # L110: 314.66666 -> 366.5
$ws.Cells.Item(110, 13).Value = 845.5  # M110: 1313.0833 -> 845.5
$ws.Cells.Item(110, 14).Value = -4456.5  # N110: -4404.66666 -> -4456.5

$ws.Cells.Item(132, 8).Value = 1900.9  # H132: 1987.1724 -> 1900.9
$ws.Cells.Item(132, 9).Value = 957.7826  # I132: 1028.6364 -> 957.7826
$ws.Cells.Item(132, 11).Value = 2873.3478  # K132: 3085.9092 -> 2873.3478
$ws.Cells.Item(132, 13).Value = -343.3478  # M132: -555.9092000000001 -> -343.3478

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(9, 8).Value = 150000  # H9: 0 -> 150000
$ws.Cells.Item(9, 10).Value = 150000  # J9: 0 -> 150000
$ws.Cells.Item(9, 12).Value = 150000  # L9: 0 -> 150000
$ws.Cells.Item(9, 14).Value = -150336  # N9: None -> -150336

$ws.Cells.Item(31, 8).Value = 2429.0833  # H31: 2429.1667 -> 2429.0833
$ws.Cells.Item(31, 9).Value = 1143.625  # I31: 1143.75 -> 1143.625
$ws.Cells.Item(31, 11).Value = 1143.625  # K31: 1143.75 -> 1143.625
$ws.Cells.Item(31, 13).Value = -848.625  # M31: -848.75 -> -848.625

$ws.Cells.Item(34, 8).Value = 2429.0833  # H34: 2429.1667 -> 2429.0833
$ws.Cells.Item(34, 9).Value = 1143.625  # I34: 1143.75 -> 1143.625
$ws.Cells.Item(34, 11).Value = 1143.625  # K34: 1143.75 -> 1143.625
$ws.Cells.Item(34, 13).Value = -941.625  # M34: -941.75 -> -941.625

$ws.Cells.Item(59, 8).Value = 37278.75  # H59: 38996.25 -> 37278.75
$ws.Cells.Item(59, 10).Value = 37278.75  # J59: 38996.25 -> 37278.75
$ws.Cells.Item(59, 12).Value = 37278.75  # L59: 38996.25 -> 37278.75
$ws.Cells.Item(59, 14).Value = -39568.75  # N59: -41286.25 -> -39568.75

$ws.Cells.Item(99, 8).Value = 2400  # H99: 2533.3333 -> 2400
$ws.Cells.Item(99, 10).Value = 0  # J99: 2800 -> 0
$ws.Cells.Item(99, 12).Value = 0  # L99: 2800 -> 0
$ws.Cells.Item(99, 14).Value = ""  # N99: clear (was -5796)

$ws.Cells.Item(126, 8).Value = 2400  # H126: 2533.3333 -> 2400
$ws.Cells.Item(126, 10).Value = 0  # J126: 2800 -> 0
$ws.Cells.Item(126, 12).Value = 0  # L126: 8400 -> 0
$ws.Cells.Item(126, 14).Value = ""  # N126: clear (was -13340)

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 99.81  # H2: 99.70999999999999 -> 99.81
$ws.Cells.Item(2, 10).Value = 99.91836499999999  # J2: 99.81632999999999 -> 99.91836499999999
$ws.Cells.Item(2, 12).Value = 599.51019  # L2: 598.89798 -> 599.51019
$ws.Cells.Item(2, 14).Value = -825.51019  # N2: -824.89798 -> -825.51019

$ws.Cells.Item(7, 8).Value = 100  # H7: 733.3333 -> 100
$ws.Cells.Item(7, 10).Value = 100  # J7: 1050 -> 100
$ws.Cells.Item(7, 12).Value = 300  # L7: 3150 -> 300
$ws.Cells.Item(7, 14).Value = -524  # N7: -3374 -> -524

$ws.Cells.Item(38, 8).Value = 707.4  # H38: 755.0714 -> 707.4
$ws.Cells.Item(38, 9).Value = 707.4  # I38: 755.0714 -> 707.4
$ws.Cells.Item(38, 11).Value = 2122.2  # K38: 2265.2142 -> 2122.2
$ws.Cells.Item(38, 13).Value = -1775.2  # M38: -1918.2142 -> -1775.2

$ws.Cells.Item(113, 8).Value = 1162.5  # H113: 1234.8 -> 1162.5
$ws.Cells.Item(113, 9).Value = 1268.75  # I113: 1424.6666 -> 1268.75
$ws.Cells.Item(113, 11).Value = 3806.25  # K113: 4273.9998 -> 3806.25
$ws.Cells.Item(113, 13).Value = -1636.25  # M113: -2103.9998 -> -1636.25

$ws.Cells.Item(128, 8).Value = 144944  # H128: 129962.664 -> 144944
$ws.Cells.Item(128, 9).Value = 144944  # I128: 129962.664 -> 144944
$ws.Cells.Item(128, 11).Value = 434832  # K128: 389887.992 -> 434832
$ws.Cells.Item(128, 13).Value = -429852  # M128: -384907.992 -> -429852

$ws.Cells.Item(131, 8).Value = 2788.1292  # H131: 2776.9375 -> 2788.1292
$ws.Cells.Item(131, 9).Value = 990  # I131: 1470 -> 990
$ws.Cells.Item(131, 11).Value = 2970  # K131: 4410 -> 2970
$ws.Cells.Item(131, 13).Value = 2070  # M131: 630 -> 2070

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 1334.3334  # H107: 1500 -> 1334.3334
$ws.Cells.Item(107, 10).Value = 1251.5  # J107: 1500 -> 1251.5
$ws.Cells.Item(107, 12).Value = 1251.5  # L107: 1500 -> 1251.5
$ws.Cells.Item(107, 14).Value = -5091.5  # N107: -5340 -> -5091.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1149.1666  # H16: 1055.7142 -> 1149.1666
$ws.Cells.Item(16, 9).Value = 1149.1666  # I16: 1055.7142 -> 1149.1666
$ws.Cells.Item(16, 11).Value = 1149.1666  # K16: 1055.7142 -> 1149.1666
$ws.Cells.Item(16, 13).Value = -979.1666  # M16: -885.7141999999999 -> -979.1666

$ws.Cells.Item(35, 8).Value = 10027.5  # H35: 5497.25 -> 10027.5
$ws.Cells.Item(35, 9).Value = 0  # I35: 999 -> 0
$ws.Cells.Item(35, 10).Value = 10027.5  # J35: 6996.6665 -> 10027.5
$ws.Cells.Item(35, 11).Value = 0  # K35: 999 -> 0
$ws.Cells.Item(35, 12).Value = 10027.5  # L35: 6996.6665 -> 10027.5
$ws.Cells.Item(35, 13).Value = ""  # M35: clear (was -663)
$ws.Cells.Item(35, 14).Value = -10699.5  # N35: -7668.6665 -> -10699.5

$ws.Cells.Item(68, 8).Value = 1866  # H68: 2933.6667 -> 1866
$ws.Cells.Item(68, 9).Value = 1866  # I68: 1899 -> 1866
$ws.Cells.Item(68, 10).Value = 0  # J68: 5003 -> 0
$ws.Cells.Item(68, 11).Value = 1866  # K68: 1899 -> 1866
$ws.Cells.Item(68, 12).Value = 0  # L68: 5003 -> 0
$ws.Cells.Item(68, 13).Value = -1117  # M68: -1150 -> -1117
$ws.Cells.Item(68, 14).Value = ""  # N68: clear (was -6501)

$ws.Cells.Item(71, 8).Value = 1866  # H71: 2933.6667 -> 1866
$ws.Cells.Item(71, 9).Value = 1866  # I71: 1899 -> 1866
$ws.Cells.Item(71, 10).Value = 0  # J71: 5003 -> 0
$ws.Cells.Item(71, 11).Value = 9330  # K71: 9495 -> 9330
$ws.Cells.Item(71, 12).Value = 0  # L71: 25015 -> 0
$ws.Cells.Item(71, 13).Value = -5586  # M71: -5751 -> -5586
$ws.Cells.Item(71, 14).Value = ""  # N71: clear (was -32503)

$ws.Cells.Item(88, 8).Value = 0  # H88: 20000 -> 0
$ws.Cells.Item(88, 10).Value = 0  # J88: 20000 -> 0
$ws.Cells.Item(88, 12).Value = 0  # L88: 20000 -> 0
$ws.Cells.Item(88, 14).Value = ""  # N88: clear (was -20856)

$ws.Cells.Item(91, 8).Value = 0  # H91: 20000 -> 0
$ws.Cells.Item(91, 10).Value = 0  # J91: 20000 -> 0
$ws.Cells.Item(91, 12).Value = 0  # L91: 20000 -> 0
$ws.Cells.Item(91, 14).Value = ""  # N91: clear (was -22964)

$ws.Cells.Item(138, 8).Value = 0  # H138: 56531.668 -> 0
$ws.Cells.Item(138, 10).Value = 0  # J138: 56531.668 -> 0
$ws.Cells.Item(138, 12).Value = 0  # L138: 56531.668 -> 0
$ws.Cells.Item(138, 14).Value = ""  # N138: clear (was -66811.66800000001)

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 20882.889  # H81: 24171.285 -> 20882.889
$ws.Cells.Item(81, 9).Value = 20882.889  # I81: 24171.285 -> 20882.889
$ws.Cells.Item(81, 11).Value = 41765.778  # K81: 48342.57 -> 41765.778
$ws.Cells.Item(81, 13).Value = -40704.778  # M81: -47281.57 -> -40704.778

$ws.Cells.Item(84, 8).Value = 20882.889  # H84: 24171.285 -> 20882.889
$ws.Cells.Item(84, 9).Value = 20882.889  # I84: 24171.285 -> 20882.889
$ws.Cells.Item(84, 11).Value = 208828.89  # K84: 241712.85 -> 208828.89
$ws.Cells.Item(84, 13).Value = -203524.89  # M84: -236408.85 -> -203524.89
